$d = $word.ActiveDocument

# --- 1) "RoleplayTalk Studios presents" -> "Red Stone Forge presents" ---
# The first (PreTitle) paragraph reads "RoleplayTalk" + " " + "Studios" + " presents"
# across several runs. Replace the "RoleplayTalk Studios " portion (leaving the
# existing " presents" run's leading space behind, since it becomes part of the
# replacement) with "Red Stone Forge ", then drop a "_GoBack" bookmark right at
# the resulting edit point (this is what Word does automatically after any
# in-place text edit), which also renumbers every other bookmark in the
# document by +1 (exactly matching the target revision).
$firstPara = $d.Paragraphs(1)
$firstRange = $firstPara.Range
$oldPhrase = "RoleplayTalk Studios "
$newPhrase = "Red Stone Forge "
$target = $d.Range($firstRange.Start, $firstRange.Start + $oldPhrase.Length)
if ($target.Text -eq $oldPhrase) {
    $target.Text = $newPhrase
}

$editPoint = $firstRange.Start + $newPhrase.Length
$goBackRange = $d.Range($editPoint, $editPoint)
$d.Bookmarks.Add("_GoBack", $goBackRange)

# --- 2) Title "Space Invaders" -> "Alien Invasion" ---
# Only the document Title (there's exactly one paragraph using the "Title"
# style) changes; the other two body mentions of "Space Invaders" stay as-is.
foreach ($p in $d.Paragraphs) {
    $style = $p.Style
    if ($style.NameLocal -eq "Title") {
        $p.Range.Find.Execute("Space Invaders", $true, $false, $false, $false, $false, `
            $true, 1, $false, "Alien Invasion", 2)
    }
}
